# Applies the Betfair Back/Lay odds update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "N3"  = 4.8
    "P3"  = 2.32
    "AB3" = 12
    "AC3" = 10
    "AF3" = 14
    "AG3" = 11
    "AN3" = 10

    "F6"  = 12
    "G6"  = 20
    "H6"  = 1.19
    "I6"  = 1.3
    "J6"  = 7
    "K6"  = 8.6

    "F8"  = 2.3
    "G8"  = 2.32
    "H8"  = 3.6
    "V8"  = 1.38
    "Z8"  = 23

    "M9"  = 1.04
    "P9"  = 2.46
    "Q9"  = 1.67
    "R9"  = 1.55
    "T9"  = 1.79
    "X9"  = 22
    "AA9" = 160
    "AG9" = 9.199999999999999
    "AM9" = 90

    "G10" = 2.12
    "J10" = 3.55
    "R10" = 1.48
    "S10" = 1.9
    "T10" = 1.04
    "U10" = 1.04

    "L11" = 1.4
    "P11" = 1.94
    "T11" = 1.81
    "U11" = 2.16
    "V11" = 1.42
    "W11" = 1.76
    "AA11" = 60
    "AC11" = 7.6
    "AD11" = 14
    "AF11" = 14
    "AG11" = 11
    "AH11" = 17.5
    "AI11" = 50
    "AJ11" = 30
    "AK11" = 25
    "AM11" = 90
    "AN11" = 19

    "F12" = 8.199999999999999
    "L12" = 1.31
    "V12" = 3.1
    "W12" = 1.13
    "Y12" = 9.6
    "Z12" = 8.800000000000001
    "AA12" = 12.5
    "AB12" = 29
    "AE12" = 14
    "AF12" = 70
    "AG12" = 29
    "AI12" = 32
    "AJ12" = 240
    "AK12" = 110
    "AN12" = 120
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
